$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.99999999612595447
$ws.Range("A2").Value = 0.99581439375090408
$ws.Range("A3").Value = 0.98111061262424404
$ws.Range("A4").Value = 0.97610381438828142
$ws.Range("A5").Value = 0.96816075845469185
$ws.Range("A6").Value = 0.94911194946342781
$ws.Range("A7").Value = 0.9435246919840159
$ws.Range("A8").Value = 0.93659926968102114
$ws.Range("A9").Value = 0.92831053368331184
$ws.Range("A10").Value = 0.92094151455204409
$ws.Range("A11").Value = 0.91994035621268133
$ws.Range("A12").Value = 0.91829154899320209
$ws.Range("A13").Value = 0.91188158326218938
$ws.Range("A14").Value = 0.90969720136275911
$ws.Range("A15").Value = 0.90710563620870843
$ws.Range("A16").Value = 0.9045990921789745
$ws.Range("A17").Value = 0.90089109259890698
$ws.Range("A18").Value = 0.89978216305739001
$ws.Range("A19").Value = 0.99312145501120108
$ws.Range("A20").Value = 0.98600426436831556
$ws.Range("A21").Value = 0.98460575369221659
$ws.Range("A22").Value = 0.98334124415017476
$ws.Range("A23").Value = 0.9740755267285699
$ws.Range("A24").Value = 0.96105436359181473
$ws.Range("A25").Value = 0.95459731426051331
$ws.Range("A26").Value = 0.94808463619151229
$ws.Range("A27").Value = 0.94690539860109224
$ws.Range("A28").Value = 0.94223938165782228
$ws.Range("A29").Value = 0.93950727168947168
$ws.Range("A30").Value = 0.93930108172330651
$ws.Range("A31").Value = 0.94762814404098106
$ws.Range("A32").Value = 0.95149151742844307
$ws.Range("A33").Value = 0.95815300697014461
